$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "-"
$ws.Range("F6").Value = "MEC-1B-Gestao Intregrada"
$ws.Range("C7").Value = "MCT-1A-Gestão integrada"
$ws.Range("F7").Value = "MEC-1B-Gestao Intregrada"
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "MCT-1A-Gestão integrada"
$ws.Range("F8").Value = "-"
$ws.Range("B11").Value = "-"
$ws.Range("B12").Value = "-"
$ws.Range("C14").Value = "MEC-1A-Gestao Integrada"
$ws.Range("C15").Value = "MEC-1A-Gestao Integrada"
$ws.Range("C20").Value = "-"
